# Update "想去人数" (want-to-go count) values in both the "展览" and
# "全部类型" sheets, which contain duplicate data tables.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 542
    $ws.Range("F4").Value = 271
    $ws.Range("F5").Value = 5
    $ws.Range("F7").Value = 775
}
